# Update the lattice-multiplication worksheet table: replace all 15 exercise
# cells with the new set of problems (per the commit's regenerated output).
#
# Each cell is a single run containing 5 text fragments separated by line
# breaks:
#   "<A> x <B>"
#   "  <B1>    <B2>"   (the two digits of the second factor)
#   "  ----"
#   "<A1>|    |"       (first digit of the first factor)
#   "<A2>|    |"       (second digit of the first factor)
#
# We rebuild each cell's Range.Text using a vertical-tab (chr 11), which
# Word's text model round-trips as a <w:br/> between runs of text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# New problems, in row-major order (5 rows x 3 columns), as
# (label, digitsRow, row1, row2) tuples.
$problems = @(
    @("59 x 68", "  6    8", "5|    |", "9|    |"),
    @("66 x 70", "  7    0", "6|    |", "6|    |"),
    @("37 x 34", "  3    4", "3|    |", "7|    |"),
    @("13 x 70", "  7    0", "1|    |", "3|    |"),
    @("59 x 68", "  6    8", "5|    |", "9|    |"),
    @("87 x 29", "  2    9", "8|    |", "7|    |"),
    @("95 x 64", "  6    4", "9|    |", "5|    |"),
    @("45 x 92", "  9    2", "4|    |", "5|    |"),
    @("35 x 19", "  1    9", "3|    |", "5|    |"),
    @("99 x 65", "  6    5", "9|    |", "9|    |"),
    @("29 x 91", "  9    1", "2|    |", "9|    |"),
    @("49 x 56", "  5    6", "4|    |", "9|    |"),
    @("73 x 21", "  2    1", "7|    |", "3|    |"),
    @("52 x 15", "  1    5", "5|    |", "2|    |"),
    @("23 x 25", "  2    5", "2|    |", "3|    |")
)

$rows = 5
$cols = 3
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $p = $problems[$i]
        $label = $p[0]
        $digits = $p[1]
        $row1 = $p[2]
        $row2 = $p[3]
        $cell = $t.Rows.Item($r).Cells.Item($c)
        $cell.Range.Text = $label + $nl + $digits + $nl + "  ----" + $nl + $row1 + $nl + $row2
        $i++
    }
}
